$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add older catch limits (Landings + Target TAC formulas) ---

# Cod / GB (row 13): Landings 12121, Target TAC = Landings / 0.16
$ws.Range("B13").Value = 12121
$ws.Range("D13").Formula = "=B13/0.16"

# Haddock / GOM (row 15): Landings 1070, Target TAC = Landings / 0.01
$ws.Range("B15").Value = 1070
$ws.Range("D15").Formula = "=B15/0.01"

# Yellowtail / CC-GOM (row 17): Target TAC = Landings / 0.09 (Landings already present)
$ws.Range("D17").Formula = "=B17/0.09"

# Pollock (row 19): Target TAC = Landings / 0.05 (Landings already present)
# Also needs to pick up the same number format / style used by the other
# Target TAC cells (row 19 previously used a different, unformatted style).
$ws.Range("D17").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("D19").Formula = "=B19/0.05"

# White Hake (row 25): Target TAC = Landings / 0.08 (Landings already present)
$ws.Range("D25").Formula = "=B25/0.08"

# --- View state changes recorded when the workbook was last saved ---
$excel.ActiveWindow.Zoom = 130
$ws.Range("E5").Select()

$wb.Application.CutCopyMode = $false
